# Update "Pais" sheet with the latest Covid-19 country data snapshot
# (countries & provincias Spain update, per commit message).
#
# Several countries changed rank order (their running totals crossed over),
# so the row that used to hold one country's data now holds another's; the
# row's position (rank) is unchanged but both the country name and the
# B:H metrics (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) are refreshed to match the new
# snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos (no rank change, values updated)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 6514376
$ws.Cells.Item(4, 3).Value = 145
$ws.Cells.Item(4, 4).Value = 3797173
$ws.Cells.Item(4, 5).Value = 2523166
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(4, 8).Value = 194037

# Rows 60/61 - Ghana & Armenia swapped rank
$ws.Cells.Item(60, 1).Value = "Armenia"
$ws.Cells.Item(60, 2).Value = 45152
$ws.Cells.Item(60, 3).Value = 199
$ws.Cells.Item(60, 4).Value = 41023
$ws.Cells.Item(60, 5).Value = 3224
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 905

$ws.Cells.Item(61, 1).Value = "Ghana"
$ws.Cells.Item(61, 2).Value = 45012
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 43898
$ws.Cells.Item(61, 5).Value = 831
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 283

# Row 66 - Afganistan (no rank change, values updated)
$ws.Cells.Item(66, 1).Value = "Afganistan"
$ws.Cells.Item(66, 2).Value = 38544
$ws.Cells.Item(66, 3).Value = 24
$ws.Cells.Item(66, 4).Value = 31048
$ws.Cells.Item(66, 5).Value = 6076
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = 1420

# Row 74 - El Salvador (no rank change, values updated)
$ws.Cells.Item(74, 1).Value = "El Salvador"
$ws.Cells.Item(74, 2).Value = 26511
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 4).Value = 16487
$ws.Cells.Item(74, 5).Value = 9254
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 5
$ws.Cells.Item(74, 8).Value = 770

# Row 75 - Australia (no rank change, values updated)
$ws.Cells.Item(75, 1).Value = "Australia"
$ws.Cells.Item(75, 2).Value = 26465
$ws.Cells.Item(75, 3).Value = 91
$ws.Cells.Item(75, 4).Value = 22863
$ws.Cells.Item(75, 5).Value = 2821
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 11
$ws.Cells.Item(75, 8).Value = 781

# Rows 96/97/98 - Malasia, Guayana Francesa & Hungria shifted rank
$ws.Cells.Item(96, 1).Value = "Hungria"
$ws.Cells.Item(96, 2).Value = 9715
$ws.Cells.Item(96, 3).Value = 411
$ws.Cells.Item(96, 4).Value = 3984
$ws.Cells.Item(96, 5).Value = 5103
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 628

$ws.Cells.Item(97, 1).Value = "Malasia"
$ws.Cells.Item(97, 2).Value = 9559
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 9136
$ws.Cells.Item(97, 5).Value = 295
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 128

$ws.Cells.Item(98, 1).Value = "Guayana Francesa"
$ws.Cells.Item(98, 2).Value = 9387
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 8946
$ws.Cells.Item(98, 5).Value = 379
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 62

# Rows 132/133 - Sri Lanka & Lituania swapped rank
$ws.Cells.Item(132, 1).Value = "Lituania"
$ws.Cells.Item(132, 2).Value = 3163
$ws.Cells.Item(132, 3).Value = 32
$ws.Cells.Item(132, 4).Value = 2008
$ws.Cells.Item(132, 5).Value = 1069
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 86

$ws.Cells.Item(133, 1).Value = "Sri Lanka"
$ws.Cells.Item(133, 2).Value = 3140
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 2935
$ws.Cells.Item(133, 5).Value = 193
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 12

# Rows 136/137 - Aruba & Bahamas swapped rank
$ws.Cells.Item(136, 1).Value = "Bahamas"
$ws.Cells.Item(136, 2).Value = 2657
$ws.Cells.Item(136, 3).Value = 72
$ws.Cells.Item(136, 4).Value = 1088
$ws.Cells.Item(136, 5).Value = 1506
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 4
$ws.Cells.Item(136, 8).Value = 63

$ws.Cells.Item(137, 1).Value = "Aruba"
$ws.Cells.Item(137, 2).Value = 2589
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 1293
$ws.Cells.Item(137, 5).Value = 1281
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 15

# Row 152 - Georgia (no rank change, values updated)
$ws.Cells.Item(152, 1).Value = "Georgia"
$ws.Cells.Item(152, 2).Value = 1773
$ws.Cells.Item(152, 3).Value = 44
$ws.Cells.Item(152, 4).Value = 1325
$ws.Cells.Item(152, 5).Value = 429
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 19

# Title - timestamp refreshed
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 09:20"
